$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("G2").Value = 70

    $ws.Range("F3").Value = 8284
    $ws.Range("G3").Value = 75

    $ws.Range("F4").Value = 144
    $ws.Range("G4").Value = "不可售"

    $ws.Range("G5").Value = "不可售"

    $ws.Range("F8").Value = 149

    $ws.Range("F10").Value = 208

    $ws.Range("F12").Value = 753

    $ws.Range("F14").Value = 5360
    $ws.Range("G14").Value = 65

    $ws.Range("F15").Value = 4

    $ws.Range("F17").Value = 86

    $ws.Range("F22").Value = 174

    $ws.Range("F23").Value = 9
}

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 1

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F24").Value = 1

$wb.Save()
